$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update speedPos (F) and zoom (G) values for rows 3-9
foreach ($r in 3..9) {
    $ws.Cells.Item($r, 6).Value = 10   # column F - speedPos
    $ws.Cells.Item($r, 7).Value = 7    # column G - zoom
}

# Update selection to F3:F9 with active cell F3
$ws.Range("F3:F9").Select()
